$d = $word.ActiveDocument

# Helper: Paragraph.Range.Text always carries a trailing paragraph-mark
# character (CR, code 13) - strip it before doing text comparisons.
function Trim-ParaText($text) {
    return $text.TrimEnd([char]13)
}

# Locate (by content) the first body paragraph that uses the default
# "Normal" style (no explicit pStyle) and already has the usual leading
# empty-run / following-text-run shape. We use it purely as a staging
# anchor so that paragraphs we insert via InsertParagraphBefore() pick up
# "no pStyle" formatting instead of inheriting a heading style.
function Get-NormalAnchorIndex {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ((Trim-ParaText $p.Range.Text) -eq "Get ready for some juicy fun with Big Max 77 - the fruit-themed slot machine with 5x3 reels and 10 active paylines. Now, I know what you're thinking - another fruit-themed slot? But trust me, this game is as sweet as it gets!") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the very top of the document.
# ---------------------------------------------------------------------------

$anchorIdx = Get-NormalAnchorIndex
$stagingAnchor = $d.Paragraphs.Item($anchorIdx)
$stagingAnchor.Range.InsertParagraphBefore()
$metaPara = $d.Paragraphs.Item($anchorIdx)

$metaText = "Meta description: Read our review of Big Max 77, a fruit-themed slot machine with high volatility, mobile compatibility, and a 97.01% RTP. Play for free and enjoy the minimalist design."
$metaPara.Range.Text = $metaText

# Bold just the "Meta description" label (first 16 characters, no colon).
$boldRange = $metaPara.Range.Duplicate
$boldRange.Start = $metaPara.Range.Start
$boldRange.End = $metaPara.Range.Start + 16
$boldRange.Bold = 1

# Move the fully-formatted paragraph (it keeps its own paragraph mark) up
# to sit right after the document title.
$metaPara.Range.Cut()

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaSlot = $d.Paragraphs.Item(2)
$metaSlot.Style = "Normal"
$metaInsertPoint = $metaSlot.Range.Duplicate
$metaInsertPoint.Collapse(1)
$metaInsertPoint.Paste()

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph that used to sit near the
#    end of the document (it was moved to the top in step 1 above).
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($i -ne 1 -and (Trim-ParaText $para.Range.Text) -eq "Play Big Max 77 for Free - High RTP and Mobile Compatibility") {
        $para.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new image
#    prompt, while keeping its italic formatting and leading empty run.
# ---------------------------------------------------------------------------

$newClosingText = 'Create a feature image for Big Max 77 to use on social media and marketing materials. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be holding a giant fruit, such as a watermelon or pineapple, and there should be a slot machine in the background with the name "Big Max 77" displayed prominently. The overall style should be fun and engaging, with bright colors and playful designs to attract potential players.'

# Build the replacement paragraph in a staging location so straight
# quotes are preserved and the italic formatting is applied cleanly.
$anchorIdx2 = Get-NormalAnchorIndex
$stagingAnchor2 = $d.Paragraphs.Item($anchorIdx2)
$stagingAnchor2.Range.InsertParagraphBefore()
$closingPara = $d.Paragraphs.Item($anchorIdx2)
$closingPara.Range.Text = $newClosingText
$closingPara.Range.Italic = 1
$closingPara.Range.Cut()

# Remove the old closing paragraph (last paragraph of the document).
$count2 = $d.Paragraphs.Count
$oldClosing = $d.Paragraphs.Item($count2)
$oldClosing.Range.Delete()

# Re-create the final paragraph slot and paste the new content into it.
$count3 = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count3)
$lastPara.Range.InsertParagraphAfter()
$closingSlot = $d.Paragraphs.Item($count3 + 1)
$closingSlot.Style = "Normal"
$closingInsertPoint = $closingSlot.Range.Duplicate
$closingInsertPoint.Collapse(1)
$closingInsertPoint.Paste()
